# ---------------------------------------------------------------------------
# Hefei comic-convention listing refresh (update gh-pages data, commit 456a3b4)
#
# The upstream scraper re-ran and produced a feed missing the duplicated
# "合肥·首届Gumi同人展" row that had been accidentally emitted twice. This
# script removes that duplicate (old row 16 on sheet "展览" / old row 17 on
# sheet "全部类型"), shifts the remaining event rows up by one, refreshes a
# handful of "want-to-go" counters (column F) that changed between scrapes,
# and drops the now-empty trailing row so the used range shrinks by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ==================== Sheet "展览" (index 1) ====================
$ws1 = $wb.Worksheets.Item(1)

# "Want-to-go" counter refreshes on rows unaffected by the row removal below
$ws1.Cells.Item(2,6).Value = 521
$ws1.Cells.Item(4,6).Value = 1526
$ws1.Cells.Item(8,6).Value = 6292
$ws1.Cells.Item(12,6).Value = 5291
$ws1.Cells.Item(15,6).Value = 1193

# Remove the duplicated "合肥·首届Gumi同人展" entry (old row 16): pull each
# following rows B:I values up by one row, row by row.
# row 16
$ws1.Cells.Item(16,2).NumberFormat = "@"
$ws1.Cells.Item(16,2).Value = '2024-07-21'
$ws1.Cells.Item(16,3).Value = '巢湖·喵喵漫游戏动漫展'
$ws1.Cells.Item(16,4).Value = '团结东路7号 巢湖宾馆'
$ws1.Cells.Item(16,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws1.Cells.Item(16,6).Value = 0
$ws1.Cells.Item(16,7).Value = 40
$ws1.Cells.Item(16,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87875'
$ws1.Cells.Item(16,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/frveR3tO1718818996505.png'
# row 17
$ws1.Cells.Item(17,2).NumberFormat = "@"
$ws1.Cells.Item(17,2).Value = '2024-07-27'
$ws1.Cells.Item(17,3).Value = '合肥·灵能百分百ONLY2.0'
$ws1.Cells.Item(17,4).Value = '铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)'
$ws1.Cells.Item(17,5).Value = '2024.07.27 10:00-07.27 17:00'
$ws1.Cells.Item(17,6).Value = 62
$ws1.Cells.Item(17,7).Value = 85
$ws1.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87497'
$ws1.Cells.Item(17,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg'
# row 18
$ws1.Cells.Item(18,2).NumberFormat = "@"
$ws1.Cells.Item(18,2).Value = '2024-07-27'
$ws1.Cells.Item(18,3).Value = '安徽·MAX特摄only展'
$ws1.Cells.Item(18,4).Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws1.Cells.Item(18,5).Value = '2024.07.27 09:30-07.27 18:00'
$ws1.Cells.Item(18,6).Value = 365
$ws1.Cells.Item(18,7).Value = 50
$ws1.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83684'
$ws1.Cells.Item(18,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg'
# row 19
$ws1.Cells.Item(19,2).NumberFormat = "@"
$ws1.Cells.Item(19,2).Value = '2024-07-27'
$ws1.Cells.Item(19,3).Value = '庐江·夏日游嘉年华'
$ws1.Cells.Item(19,4).Value = '白山路东150米 庐江体育馆'
$ws1.Cells.Item(19,5).Value = '2024.07.27 09:00-07.28 17:00'
$ws1.Cells.Item(19,6).Value = 72
$ws1.Cells.Item(19,7).Value = 60
$ws1.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87569'
$ws1.Cells.Item(19,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg'
# row 20
$ws1.Cells.Item(20,2).NumberFormat = "@"
$ws1.Cells.Item(20,2).Value = '2024-07-27'
$ws1.Cells.Item(20,3).Value = '长丰·莓可可游戏动漫展'
$ws1.Cells.Item(20,4).Value = '长寿路12号 长丰宾馆·梅山饭店(长寿路店)'
$ws1.Cells.Item(20,5).Value = '2024.07.27 10:00-07.27 17:00'
$ws1.Cells.Item(20,6).Value = 12
$ws1.Cells.Item(20,7).Value = 40
$ws1.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87796'
$ws1.Cells.Item(20,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/MLTfeikq1718823574810.png'
# row 21
$ws1.Cells.Item(21,2).NumberFormat = "@"
$ws1.Cells.Item(21,2).Value = '2024-07-28'
$ws1.Cells.Item(21,3).Value = '合肥·咒术回战only'
$ws1.Cells.Item(21,4).Value = '清河路19号 依立腾工业园区'
$ws1.Cells.Item(21,5).Value = '2024.07.28 09:30-07.28 17:30'
$ws1.Cells.Item(21,6).Value = 305
$ws1.Cells.Item(21,7).Value = 60
$ws1.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86520'
$ws1.Cells.Item(21,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png'
# row 22
$ws1.Cells.Item(22,2).NumberFormat = "@"
$ws1.Cells.Item(22,2).Value = '2024-07-28'
$ws1.Cells.Item(22,3).Value = '合肥·第二届TH动漫游戏嘉年华'
$ws1.Cells.Item(22,4).Value = '田埠西路199号 吉祥如意宴会楼蜀山店'
$ws1.Cells.Item(22,5).Value = '2024.07.28 09:30-07.28 17:00'
$ws1.Cells.Item(22,6).Value = 29
$ws1.Cells.Item(22,7).Value = 55
$ws1.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87447'
$ws1.Cells.Item(22,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png'
# row 23
$ws1.Cells.Item(23,2).NumberFormat = "@"
$ws1.Cells.Item(23,2).Value = '2024-08-03'
$ws1.Cells.Item(23,3).Value = '合肥·第七届环形宇宙动漫游戏嘉年华'
$ws1.Cells.Item(23,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws1.Cells.Item(23,5).Value = '2024.08.03 09:30-08.04 17:00'
$ws1.Cells.Item(23,6).Value = 3780
$ws1.Cells.Item(23,7).Value = 49
$ws1.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84767'
$ws1.Cells.Item(23,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'
# row 24
$ws1.Cells.Item(24,2).NumberFormat = "@"
$ws1.Cells.Item(24,2).Value = '2024-08-17'
$ws1.Cells.Item(24,3).Value = '合肥·银魂主题派对only2.0'
$ws1.Cells.Item(24,4).Value = '长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)'
$ws1.Cells.Item(24,5).Value = '2024.08.17 13:00-08.17 18:00'
$ws1.Cells.Item(24,6).Value = 163
$ws1.Cells.Item(24,7).Value = 128
$ws1.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87173'
$ws1.Cells.Item(24,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png'

# The data that used to live in the final row has now been copied into the
# row above it; delete the trailing row so the sheets used range is A1:I24.
$ws1.Rows.Item(25).Delete()

# ==================== Sheet "全部类型" (index 4) ====================
$ws4 = $wb.Worksheets.Item(4)

# "Want-to-go" counter refreshes on rows unaffected by the row removal below
$ws4.Cells.Item(3,6).Value = 521
$ws4.Cells.Item(5,6).Value = 1526
$ws4.Cells.Item(9,6).Value = 6292
$ws4.Cells.Item(13,6).Value = 5291
$ws4.Cells.Item(16,6).Value = 1193

# Remove the duplicated "合肥·首届Gumi同人展" entry (old row 17): pull each
# following rows B:I values up by one row, row by row.
# row 17
$ws4.Cells.Item(17,2).NumberFormat = "@"
$ws4.Cells.Item(17,2).Value = '2024-07-21'
$ws4.Cells.Item(17,3).Value = '巢湖·喵喵漫游戏动漫展'
$ws4.Cells.Item(17,4).Value = '团结东路7号 巢湖宾馆'
$ws4.Cells.Item(17,5).Value = '2024.07.21 10:00-07.21 17:00'
$ws4.Cells.Item(17,6).Value = 0
$ws4.Cells.Item(17,7).Value = 40
$ws4.Cells.Item(17,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87875'
$ws4.Cells.Item(17,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/frveR3tO1718818996505.png'
# row 18
$ws4.Cells.Item(18,2).NumberFormat = "@"
$ws4.Cells.Item(18,2).Value = '2024-07-27'
$ws4.Cells.Item(18,3).Value = '合肥·灵能百分百ONLY2.0'
$ws4.Cells.Item(18,4).Value = '铜陵北路金邦国际大厦一楼 格律诗婚礼艺术中心(新站店)'
$ws4.Cells.Item(18,5).Value = '2024.07.27 10:00-07.27 17:00'
$ws4.Cells.Item(18,6).Value = 62
$ws4.Cells.Item(18,7).Value = 85
$ws4.Cells.Item(18,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87497'
$ws4.Cells.Item(18,9).Value = '//i1.hdslb.com/bfs/openplatform/202406/3Jycwu1U1717858639976.jpeg'
# row 19
$ws4.Cells.Item(19,2).NumberFormat = "@"
$ws4.Cells.Item(19,2).Value = '2024-07-27'
$ws4.Cells.Item(19,3).Value = '安徽·MAX特摄only展'
$ws4.Cells.Item(19,4).Value = '桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间'
$ws4.Cells.Item(19,5).Value = '2024.07.27 09:30-07.27 18:00'
$ws4.Cells.Item(19,6).Value = 365
$ws4.Cells.Item(19,7).Value = 50
$ws4.Cells.Item(19,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83684'
$ws4.Cells.Item(19,9).Value = '//i0.hdslb.com/bfs/openplatform/202405/qBnW1VeB1715423018997.jpeg'
# row 20
$ws4.Cells.Item(20,2).NumberFormat = "@"
$ws4.Cells.Item(20,2).Value = '2024-07-27'
$ws4.Cells.Item(20,3).Value = '庐江·夏日游嘉年华'
$ws4.Cells.Item(20,4).Value = '白山路东150米 庐江体育馆'
$ws4.Cells.Item(20,5).Value = '2024.07.27 09:00-07.28 17:00'
$ws4.Cells.Item(20,6).Value = 72
$ws4.Cells.Item(20,7).Value = 60
$ws4.Cells.Item(20,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87569'
$ws4.Cells.Item(20,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/5tB3RWrN1718243791381.jpeg'
# row 21
$ws4.Cells.Item(21,2).NumberFormat = "@"
$ws4.Cells.Item(21,2).Value = '2024-07-27'
$ws4.Cells.Item(21,3).Value = '长丰·莓可可游戏动漫展'
$ws4.Cells.Item(21,4).Value = '长寿路12号 长丰宾馆·梅山饭店(长寿路店)'
$ws4.Cells.Item(21,5).Value = '2024.07.27 10:00-07.27 17:00'
$ws4.Cells.Item(21,6).Value = 12
$ws4.Cells.Item(21,7).Value = 40
$ws4.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87796'
$ws4.Cells.Item(21,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/MLTfeikq1718823574810.png'
# row 22
$ws4.Cells.Item(22,2).NumberFormat = "@"
$ws4.Cells.Item(22,2).Value = '2024-07-28'
$ws4.Cells.Item(22,3).Value = '合肥·咒术回战only'
$ws4.Cells.Item(22,4).Value = '清河路19号 依立腾工业园区'
$ws4.Cells.Item(22,5).Value = '2024.07.28 09:30-07.28 17:30'
$ws4.Cells.Item(22,6).Value = 305
$ws4.Cells.Item(22,7).Value = 60
$ws4.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=86520'
$ws4.Cells.Item(22,9).Value = '//i2.hdslb.com/bfs/openplatform/202405/cLCM0a1e1716952386781.png'
# row 23
$ws4.Cells.Item(23,2).NumberFormat = "@"
$ws4.Cells.Item(23,2).Value = '2024-07-28'
$ws4.Cells.Item(23,3).Value = '合肥·第二届TH动漫游戏嘉年华'
$ws4.Cells.Item(23,4).Value = '田埠西路199号 吉祥如意宴会楼蜀山店'
$ws4.Cells.Item(23,5).Value = '2024.07.28 09:30-07.28 17:00'
$ws4.Cells.Item(23,6).Value = 29
$ws4.Cells.Item(23,7).Value = 55
$ws4.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87447'
$ws4.Cells.Item(23,9).Value = '//i0.hdslb.com/bfs/openplatform/202406/jHqfdzLQ1718091324240.png'
# row 24
$ws4.Cells.Item(24,2).NumberFormat = "@"
$ws4.Cells.Item(24,2).Value = '2024-08-03'
$ws4.Cells.Item(24,3).Value = '合肥·第七届环形宇宙动漫游戏嘉年华'
$ws4.Cells.Item(24,4).Value = '南京路与庐州大道交汇处 合肥滨湖国际会展中心'
$ws4.Cells.Item(24,5).Value = '2024.08.03 09:30-08.04 17:00'
$ws4.Cells.Item(24,6).Value = 3780
$ws4.Cells.Item(24,7).Value = 49
$ws4.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=84767'
$ws4.Cells.Item(24,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/nBGuQecO1713856894035.jpeg'
# row 25
$ws4.Cells.Item(25,2).NumberFormat = "@"
$ws4.Cells.Item(25,2).Value = '2024-08-03'
$ws4.Cells.Item(25,3).Value = '合肥·首届包河留声机音乐节—《菊次郎的夏天》久石让钢琴曲梦幻之旅演奏会'
$ws4.Cells.Item(25,4).Value = '徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院'
$ws4.Cells.Item(25,5).Value = '2024.08.03 19:30-08.03 21:00'
$ws4.Cells.Item(25,6).Value = 42
$ws4.Cells.Item(25,7).Value = 80
$ws4.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83556'
$ws4.Cells.Item(25,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/4nwOTVDu1711695345941.jpeg'
# row 26
$ws4.Cells.Item(26,2).NumberFormat = "@"
$ws4.Cells.Item(26,2).Value = '2024-08-17'
$ws4.Cells.Item(26,3).Value = '合肥·银魂主题派对only2.0'
$ws4.Cells.Item(26,4).Value = '长江东路1137号圣大国际商贸中心2-301室 梦田音乐LiveHouse(合肥店)'
$ws4.Cells.Item(26,5).Value = '2024.08.17 13:00-08.17 18:00'
$ws4.Cells.Item(26,6).Value = 163
$ws4.Cells.Item(26,7).Value = 128
$ws4.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=87173'
$ws4.Cells.Item(26,9).Value = '//i2.hdslb.com/bfs/openplatform/202406/aSc8SoTl1718078234193.png'

# The data that used to live in the final row has now been copied into the
# row above it; delete the trailing row so the sheets used range is A1:I26.
$ws4.Rows.Item(27).Delete()

